$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# R2: empty cell, same format as Q2
$ws.Range("Q2").Copy() | Out-Null
$ws.Range("R2").PasteSpecial(-4122) | Out-Null

# R3: 2021, same format as Q3
$ws.Range("R3").Value = 2021
$ws.Range("Q3").Copy() | Out-Null
$ws.Range("R3").PasteSpecial(-4122) | Out-Null

# R4: 13.5, same format as Q4
$ws.Range("R4").Value = 13.5
$ws.Range("Q4").Copy() | Out-Null
$ws.Range("R4").PasteSpecial(-4122) | Out-Null

# R5: 15.1, same format as Q5
$ws.Range("R5").Value = 15.1
$ws.Range("Q5").Copy() | Out-Null
$ws.Range("R5").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Update the selection to T3 as in the target workbook
$ws.Range("T3").Select() | Out-Null
